$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F holds "想去人数" (want-to-go count)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 5397
$wsExhibition.Range("F4").Value = 11478
$wsExhibition.Range("F6").Value = 588
$wsExhibition.Range("F8").Value = 260
$wsExhibition.Range("F9").Value = 991

# Sheet "全部类型" (All types) - same events repeated, different rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 5397
$wsAll.Range("F7").Value = 11478
$wsAll.Range("F9").Value = 588
$wsAll.Range("F13").Value = 260
$wsAll.Range("F14").Value = 991
